$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "change proneness" values in column C (rows 2-15) with the
# newly-computed correlation data (metric1&2&6 and metric5&6 analysis).
$values = @(
    0.005917159763313608,
    0.020833333333333343,
    0.017241379310344838,
    0.03125,
    0.04166666666666668,
    0.025641025641025637,
    0.0017730496453900704,
    0.005076142131979692,
    0.0030303030303030294,
    0.009900990099009917,
    0.07142857142857142,
    0.0052356020942408415,
    0.04545454545454545,
    0.0007704160246533115
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Update the active selection to reflect where the author was last working.
$ws.Range("E11").Select()
